# ============================================================================
# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and
# excel sheets
#
# 1) Clean up a handful of cells on "ODI Batting Extra" that were being
#    written as empty placeholders - they should simply not exist.
# 2) Add a brand-new "ODI Bowling Extra" sheet (mirrors "ODI Batting Extra"
#    but for bowling stats) and populate it with the scraped data.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Strip the stray empty cells out of "ODI Batting Extra"
# ---------------------------------------------------------------------------
$wsBattingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$wsBattingExtra.Range("E3").ClearContents()
$wsBattingExtra.Range("B5:E5").ClearContents()
$wsBattingExtra.Range("C9:E9").ClearContents()
$wsBattingExtra.Range("B10:E10").ClearContents()
$wsBattingExtra.Range("C11:E11").ClearContents()
$wsBattingExtra.Range("B18:E18").ClearContents()

# ---------------------------------------------------------------------------
# 2) Create the new "ODI Bowling Extra" sheet at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsBowlingExtra = $wb.Worksheets.Add($null, $lastSheet)
$wsBowlingExtra.Name = "ODI Bowling Extra"

# Header row
$wsBowlingExtra.Range("A1").Value = "MATCH_CODE"
$wsBowlingExtra.Range("B1").Value = "MAIDEN_OVERS"
$wsBowlingExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"
$wsBowlingExtra.Range("A1:C1").Style = $wsBattingExtra.Range("A1").Style

# Data rows: MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL
$data = @(
    @("4218", "0", ""),
    @("4229", "0", ""),
    @("4230", "0", ""),
    @("4250", "0", "10.00%"),
    @("4251", "0", ""),
    @("4252", "0", ""),
    @("4356", "", ""),
    @("4358", "", ""),
    @("4416", "0", ""),
    @("4418", "", ""),
    @("4447", "0", ""),
    @("4463", "0", ""),
    @("4481", "0", "20.00%"),
    @("4537", "0", "10.00%"),
    @("4538", "", ""),
    @("4550", "0", "10.00%"),
    @("4557", "0", ""),
    @("4626", "0", ""),
    @("4682", "0", "10.00%"),
    @("4685", "", "")
)

$row = 2
foreach ($item in $data) {
    $matchCode = $item[0]
    $maidenOvers = $item[1]
    $percentWickets = $item[2]

    $cellA = $wsBowlingExtra.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $matchCode

    if ($maidenOvers -ne "") {
        $cellB = $wsBowlingExtra.Cells.Item($row, 2)
        $cellB.NumberFormat = "@"
        $cellB.Value = $maidenOvers
    }

    if ($percentWickets -ne "") {
        $cellC = $wsBowlingExtra.Cells.Item($row, 3)
        $cellC.NumberFormat = "@"
        $cellC.Value = $percentWickets
    }

    $row++
}
